$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column "type" before the existing catchment.id column (old column B).
$ws.Range("B1").EntireColumn.Insert()

# Site "type" classification values for each data row (rows 2-16). These are entered
# before the header text below so the shared-string table indices line up with the
# original authoring order (Control, Mixed, Harvest, Insect, then "type").
$ws.Range("B2").Value = "Control"
$ws.Range("B3").Value = "Control"
$ws.Range("B4").Value = "Control"
$ws.Range("B5").Value = "Mixed"
$ws.Range("B6").Value = "Control"
$ws.Range("B7").Value = "Control"
$ws.Range("B8").Value = "Harvest"
$ws.Range("B9").Value = "Harvest"
$ws.Range("B10").Value = "Control"
$ws.Range("B11").Value = "Insect"
$ws.Range("B12").Value = "Insect"
$ws.Range("B13").Value = "Mixed"
$ws.Range("B14").Value = "Insect"
$ws.Range("B15").Value = "Mixed"
$ws.Range("B16").Value = "Mixed"

# Header for the new column.
$ws.Range("B1").Value = "type"

# The column insert shifts the existing array formulas one column to the right but
# breaks their array ("Ctrl+Shift+Enter") nature, so re-enter them explicitly.
$ws.Range("E2:E16").FormulaArray = "=(D2:D16/1000)"
$ws.Range("G2:G16").FormulaArray = "=(F2:F16/1000)"
$ws.Range("I2:I16").FormulaArray = "=(H2:H16/1000)"
$ws.Range("K2:K16").FormulaArray = "=(J2:J16/1000)"

# View changes: zoom level and selected cells.
$excel.ActiveWindow.Zoom = 160

$u = $ws.Range("A14:XFD14,A11:XFD11,A12:XFD12")
$u.Select()

Write-Host "done"
